$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("#Concerns")

# Rename the sheet's header labels ([Concerns] -> [Objectives], ttIsaConcern -> ttIsaObjective)
$ws.Range("A1").Value = "[Objectives]"
$ws.Range("B1").Value = "ttIsaObjective"

# Update the "Conc_" prefix formulas to "Obj_" (column A anchor formulas)
$ws.Range("A3").Formula = '=IF(OR($C3="",$D3=""),"",CONCATENATE("Obj_",$C3,"_",$D3))'
$ws.Range("A4:A18").Formula = '=IF(OR($C4="",$D4=""),"",CONCATENATE("Obj_",$C4,"_",$D4))'

# Hide column B (duplicate label column) and update the selection accordingly
$ws.Columns("B").Hidden = $true
$ws.Range("B1:B1048576").Select()
